# The "G" column (header "K") previously held a "Strike#" style value.
# Regenerate the save_data so that column G now reports the K value
# (std/mean derived "s_vals") for each row, and write the updated figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 2
    4  = 0
    5  = 1
    6  = 1
    7  = 2
    8  = 1
    9  = 5
    10 = 3
    11 = 2
    12 = 3
    13 = 0
    14 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
